# Missing fields of RTW
#
# The "Remuneration_Start", "Remuneration_End" and "Measure_months" field
# rows (rows 3-5 on Sheet1) are no longer needed, so this removes them
# entirely. Excel automatically shifts the remaining rows up, drops the
# now-unused shared strings, and re-points every remaining cell at its
# (unchanged) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three obsolete field rows (Remuneration_Start, Remuneration_End,
# Measure_months). This shifts old rows 6-14 up to become new rows 3-11.
$ws.Rows("3:5").Delete() | Out-Null

# The row that now holds "cost_code" (row 10) inherited the custom 60pt row
# height that used to belong to the row below it before the shift; reset it
# back to the sheet's default height so only the "cost_code2" row (row 11)
# keeps its explicit 60pt height, matching the original layout intent.
$ws.Rows(10).AutoFit() | Out-Null

# Update the active cell selection to match the saved view.
$ws.Range("F8").Select() | Out-Null
